$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("摘要" / abstract) entirely - remaining columns D,E,F,G
# shift left to C,D,E,F, keeping their widths/contents intact.
$ws.Columns.Item(3).Delete()

# Rename A1 from "编号" to "0-编号"
$ws.Range("A1").Value = "0-编号"

# Column B ("题目"/title) needs to be much wider to hold full titles.
# (149.54296875 in OOXML char-width units -> closest COM ColumnWidth
# input that survives this host's pixel-grid rounding.)
$ws.Columns.Item(2).ColumnWidth = 148.85714285714286

# Center-align the header row (A1:F1)
$ws.Range("A1:F1").HorizontalAlignment = -4108

# Update selection to match target state
$ws.Range("B9").Select()
